# This workbook holds one weekly "Caqui" (persimmon) price record per row
# (rows 2-16). The update re-distributes the per-report fields (date,
# quality, volume, min/max/avg price, sale unit, origin, $/Kg, Kg/unit)
# across the existing rows - i.e. the set of records stays the same, but
# each row ends up showing a different record's data. The identifying /
# descriptive columns (A,B,C,E,F,G,H,I,J,K - market, product, variety...)
# are identical on every row already and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record values which get reshuffled.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D, L, M, N, O, P, Q, R, S, T

# Snapshot the current ("before") values for every affected row/column so the
# reassignment below can be computed from a stable source regardless of the
# order cells get overwritten in.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2()
    }
    $snapshot[$r] = $rowData
}

# Maps each destination row to the row whose (D,L,M,N,O,P,Q,R,S,T) values it
# should receive.
$mapping = @{
    2  = 4
    3  = 15
    4  = 11
    5  = 12
    6  = 14
    7  = 13
    8  = 16
    9  = 9
    10 = 6
    11 = 3
    12 = 8
    13 = 7
    14 = 10
    15 = 2
    16 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcData[$c]
    }
}
